$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.115.25'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +2.00%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.927.72'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.09%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '482.55'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +3.10%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.49'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.54%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.622'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -2.47%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.997'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.07%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.723'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -3.49%  '

# Row 10
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +8.37%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000360'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +13.53%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '42.58'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -3.20%  '

# Row 13
$ws.Range('B13').NumberFormat = '@'
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').NumberFormat = '@'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.554.35'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.16%  '

# Row 14
$ws.Range('B14').NumberFormat = '@'
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').NumberFormat = '@'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.46'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.06%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.937.63'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.66%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.56'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.80%  '

# Row 17
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.36%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '19.70'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -2.23%  '

# Row 19
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -3.73%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '69.087.31'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.53%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '434.58'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.41%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '14.61'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.61%  '

# Row 23
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.44%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '88.12'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.79%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.64'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +14.87%  '

# Row 26
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.98%  '

# Row 27
$ws.Range('B27').NumberFormat = '@'
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').NumberFormat = '@'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '38.25'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.01%  '

# Row 28
$ws.Range('B28').NumberFormat = '@'
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').NumberFormat = '@'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.39'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +2.03%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.90'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +7.94%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '712.33'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -3.21%  '

# Row 31
$ws.Range('B31').NumberFormat = '@'
$ws.Range('B31').Value = 'Cosmos'
$ws.Range('C31').NumberFormat = '@'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '13.23'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -4.26%  '

# Row 32
$ws.Range('B32').NumberFormat = '@'
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').NumberFormat = '@'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.129'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -4.76%  '

# Row 33
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +2.60%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0₃0934'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +34.72%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '41.18'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -4.80%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '58.72'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +1.55%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.152'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -7.12%  '

# Row 38
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.63%  '

# Row 39
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.23%  '

# Row 40
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.42%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.76'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +8.25%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.00'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +7.27%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.98'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +1.76%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.339'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.88%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.141'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -1.25%  '

# Row 46
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.01%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.41'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -3.01%  '

# Row 48
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.55%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '148.24'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +2.70%  '

# Row 50
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -3.34%  '

# Row 51
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -2.83%  '
